# Add a new "bacteria" data column (column S) to the FLP raw-data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1): label the new column, keeping its existing style (s="2") ---
$ws.Range("S1").Value = "bacteria"

# --- Apply the same number format/font used by the neighbouring column T
#     (style index 4: numFmt "0.00E+00", Aptos Narrow font) to all of the
#     new data cells S2:S22 before writing the values, so no duplicate
#     style entries get created. ---
$ws.Range("T2").Copy()
$ws.Range("S2:S22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data values (row -> bacteria count) ---
$ws.Range("S2").Value  = 381026.94199999998
$ws.Range("S3").Value  = 317831.799
$ws.Range("S4").Value  = 1036428.52
$ws.Range("S5").Value  = 959287.26199999999
$ws.Range("S6").Value  = 1542260.84
$ws.Range("S7").Value  = 1454528.5
$ws.Range("S8").Value  = 940495.71699999995
$ws.Range("S9").Value  = 907026.16799999995
$ws.Range("S10").Value = 371642.77600000001
$ws.Range("S11").Value = 490596.967
$ws.Range("S12").Value = 426256.78200000001
$ws.Range("S13").Value = 338137.217
$ws.Range("S14").Value = 357910.80599999998
$ws.Range("S15").Value = 1060462.81
$ws.Range("S16").Value = 990391.18299999996
$ws.Range("S17").Value = 443989.772
$ws.Range("S18").Value = 513176.005
$ws.Range("S19").Value = 912927.25399999996
$ws.Range("S20").Value = 813150.12300000002
$ws.Range("S21").Value = 522691.05900000001
$ws.Range("S22").Value = 515083.51400000002

# --- Update the view: zoom to 69% and move the active selection to U9 ---
$excel.ActiveWindow.Zoom = 69
$ws.Range("U9").Select()
